$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value from 137 to 200
$ws.Range("B2").Value = 200

# Delete row 4 entirely (A4=2, B4=63)
$ws.Rows(4).Delete()
